$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 137
$ws.Range("A137").Value = 112111486
$ws.Range("B137").Value = 83086
$ws.Range("E137").Value = 5589
$ws.Range("F137").Value = 'Rödbrun klubbdyna'
$ws.Range("G137").Value = 'Trichoderma nybergianum'
$ws.Range("H137").Value = '(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr'
$ws.Range("Q137").Value = 446833
$ws.Range("R137").Value = 7032727

# Row 138
$ws.Range("A138").Value = 112111606
$ws.Range("B138").Value = 89047
$ws.Range("E138").Value = 3286
$ws.Range("F138").Value = 'Flattoppad klubbsvamp'
$ws.Range("G138").Value = 'Clavariadelphus truncatus'
$ws.Range("H138").Value = '(Quél.) Donk'
$ws.Range("Q138").Value = 446863
$ws.Range("R138").Value = 7032718

# Row 139
$ws.Range("A139").Value = 112103325
$ws.Range("B139").Value = 85448
$ws.Range("E139").Value = 3739
$ws.Range("F139").Value = 'Persiljespindling'
$ws.Range("G139").Value = 'Cortinarius sulfurinus'
$ws.Range("H139").Value = 'Quél.'
$ws.Range("Q139").Value = 446867
$ws.Range("R139").Value = 7032726

# Row 140
$ws.Range("A140").Value = 112104547
$ws.Range("B140").Value = 89094
$ws.Range("D140").Value = 'VU'
$ws.Range("E140").Value = 256335
$ws.Range("F140").Value = 'Taggfingersvamp'
$ws.Range("G140").Value = 'Ramaria karstenii'
$ws.Range("H140").Value = '(Sacc. & P.Syd.) Corner'
$ws.Range("Q140").Value = 446688
$ws.Range("R140").Value = 7032560

# Row 141
$ws.Range("A141").Value = 112111498
$ws.Range("B141").Value = 88181
$ws.Range("E141").Value = 1599
$ws.Range("F141").Value = 'Fjällfotad musseron'
$ws.Range("G141").Value = 'Tricholoma olivaceotinctum'
$ws.Range("H141").Value = 'Mort.Chr. & Heilm.-Claus.'
$ws.Range("Q141").Value = 446860
$ws.Range("R141").Value = 7032743

# Row 142
$ws.Range("A142").Value = 112102200
$ws.Range("B142").Value = 89047
$ws.Range("E142").Value = 3286
$ws.Range("F142").Value = 'Flattoppad klubbsvamp'
$ws.Range("G142").Value = 'Clavariadelphus truncatus'
$ws.Range("H142").Value = '(Quél.) Donk'
$ws.Range("Q142").Value = 446961
$ws.Range("R142").Value = 7032566

# Row 143
$ws.Range("A143").Value = 112101773
$ws.Range("B143").Value = 86371
$ws.Range("D143").Value = 'NT'
$ws.Range("E143").Value = 4412
$ws.Range("F143").Value = 'Äggvaxskivling'
$ws.Range("G143").Value = 'Hygrophorus karstenii'
$ws.Range("H143").Value = 'Sacc. & Cub.'
$ws.Range("Q143").Value = 446984
$ws.Range("R143").Value = 7032942

# Row 144
$ws.Range("B144").Value = 89104

# Row 145
$ws.Range("A145").Value = 112104573
$ws.Range("B145").Value = 88181
$ws.Range("D145").Value = 'VU'
$ws.Range("E145").Value = 1599
$ws.Range("F145").Value = 'Fjällfotad musseron'
$ws.Range("G145").Value = 'Tricholoma olivaceotinctum'
$ws.Range("H145").Value = 'Mort.Chr. & Heilm.-Claus.'
$ws.Range("Q145").Value = 446696
$ws.Range("R145").Value = 7032530

# Row 146
$ws.Range("A146").Value = 112104270
$ws.Range("B146").Value = 85448
$ws.Range("D146").Value = 'NT'
$ws.Range("E146").Value = 3739
$ws.Range("F146").Value = 'Persiljespindling'
$ws.Range("G146").Value = 'Cortinarius sulfurinus'
$ws.Range("H146").Value = 'Quél.'
$ws.Range("Q146").Value = 446732
$ws.Range("R146").Value = 7032598

# Row 147
$ws.Range("A147").Value = 112102104
$ws.Range("B147").Value = 90799
$ws.Range("D147").Value = 'NT'
$ws.Range("E147").Value = 1968
$ws.Range("F147").Value = 'Grantaggsvamp'
$ws.Range("G147").Value = 'Bankera violascens'
$ws.Range("H147").Value = '(Alb. & Schwein. : Fr.) Pouzar'
$ws.Range("Q147").Value = 446883
$ws.Range("R147").Value = 7032568

# Row 148
$ws.Range("A148").Value = 112104553
$ws.Range("B148").Value = 86371
$ws.Range("E148").Value = 4412
$ws.Range("F148").Value = 'Äggvaxskivling'
$ws.Range("G148").Value = 'Hygrophorus karstenii'
$ws.Range("H148").Value = 'Sacc. & Cub.'
$ws.Range("Q148").Value = 446688
$ws.Range("R148").Value = 7032560

# Row 149
$ws.Range("A149").Value = 112102682
$ws.Range("B149").Value = 84955
$ws.Range("D149").Value = 'VU'
$ws.Range("E149").Value = 275
$ws.Range("F149").Value = 'Kejsarskivling'
$ws.Range("G149").Value = 'Catathelasma imperiale'
$ws.Range("H149").Value = '(P.Karst.) Singer'
$ws.Range("Q149").Value = 447025
$ws.Range("R149").Value = 7032672

# Row 150
$ws.Range("A150").Value = 112102196
$ws.Range("B150").Value = 89104
$ws.Range("D150").Value = 'VU'
$ws.Range("E150").Value = 5747
$ws.Range("F150").Value = 'Läderdoftande fingersvamp'
$ws.Range("G150").Value = 'Ramaria safraniolens'
$ws.Range("H150").Value = 'Christian'
$ws.Range("Q150").Value = 446964
$ws.Range("R150").Value = 7032565

# Row 151
$ws.Range("A151").Value = 112101944
$ws.Range("B151").Value = 89047
$ws.Range("D151").Value = 'NT'
$ws.Range("E151").Value = 3286
$ws.Range("F151").Value = 'Flattoppad klubbsvamp'
$ws.Range("G151").Value = 'Clavariadelphus truncatus'
$ws.Range("H151").Value = '(Quél.) Donk'
$ws.Range("Q151").Value = 446858
$ws.Range("R151").Value = 7032743
